# Updated with UnixDateTime Column
# Converts column E (previously a text concatenation of Date & Time) into a
# real numeric date/time value, and adds a new column P ("UnixDateTime")
# that converts that date/time into a Unix timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column E: replace the "=C&" - "&D" text formula with a literal
#    numeric date/time serial value, formatted as a date+time.
# ---------------------------------------------------------------------------
$eValues = New-Object 'object[]' 88

$eValues[2] = 41078.395833333336
$eValues[3] = 41079.395833333336
$eValues[4] = 41080.395833333336
$eValues[5] = 41081.395833333336
$eValues[6] = 41081.395833333336
$eValues[7] = 41085.395833333336
$eValues[8] = 41086.395833333336
$eValues[9] = 41087.395833333336
$eValues[10] = 41093.395833333336
$eValues[11] = 41094.395833333336
$eValues[12] = 41096.395833333336
$eValues[13] = 41099.395833333336
$eValues[14] = 41101.416666666664
$eValues[15] = 41123.75
$eValues[16] = 41144.395833333336
$eValues[17] = 41145.395833333336
$eValues[18] = 41148.395833333336
$eValues[19] = 41149.395833333336
$eValues[20] = 41149.416666666664
$eValues[21] = 41150.416666666664
$eValues[22] = 41151.395833333336
$eValues[23] = 41157.395833333336
$eValues[24] = 41158.395833333336
$eValues[25] = 41163.395833333336
$eValues[26] = 41164.395833333336
$eValues[27] = 41165.395833333336
$eValues[28] = 41169.395833333336
$eValues[29] = 41171.395833333336
$eValues[30] = 41172.395833333336
$eValues[31] = 41172.75
$eValues[32] = 41177.395833333336
$eValues[33] = 41178.395833333336
$eValues[34] = 41179.395833333336
$eValues[35] = 41184.395833333336
$eValues[36] = 41185.395833333336
$eValues[37] = 41186.791666666664
$eValues[38] = 41191.395833333336
$eValues[39] = 41192.416666666664
$eValues[40] = 41197.395833333336
$eValues[41] = 41198.395833333336
$eValues[42] = 41199.395833333336
$eValues[43] = 41200.395833333336
$eValues[44] = 41200.75
$eValues[45] = 41205.395833333336
$eValues[46] = 41206.416666666664
$eValues[47] = 41206.645833333336
$eValues[48] = 41214.395833333336
$eValues[49] = 41214.395833333336
$eValues[50] = 41219.395833333336
$eValues[51] = 41220.395833333336
$eValues[52] = 41226.395833333336
$eValues[53] = 41227.416666666664
$eValues[54] = 41228.395833333336
$eValues[55] = 41229.395833333336
$eValues[56] = 41232.395833333336
$eValues[57] = 41233.395833333336
$eValues[58] = 41234.395833333336
$eValues[59] = 41240.395833333336
$eValues[60] = 41241.416666666664
$eValues[61] = 41242.5625
$eValues[62] = 41247.395833333336
$eValues[63] = 41248.5
$eValues[64] = 41249.395833333336
$eValues[65] = 41253.541666666664
$eValues[66] = 41254.395833333336
$eValues[67] = 41255.416666666664
$eValues[68] = 41256.416666666664
$eValues[69] = 41262.416666666664
$eValues[70] = 41283.416666666664
$eValues[71] = 41288.395833333336
$eValues[72] = 41289.395833333336
$eValues[73] = 41290.395833333336
$eValues[74] = 41291.395833333336
$eValues[75] = 41291.5625
$eValues[76] = 41291.75
$eValues[77] = 41295.395833333336
$eValues[78] = 41297.416666666664
$eValues[79] = 41304.416666666664
$eValues[80] = 41310.395833333336
$eValues[81] = 41311.395833333336
$eValues[82] = 41312.395833333336
$eValues[83] = 41317.395833333336
$eValues[84] = 41318.416666666664
$eValues[85] = 41324.395833333336
$eValues[86] = 41325.395833333336
$eValues[87] = 41326.395833333336

for ($r = 2; $r -le 87; $r++) {
    $ws.Cells.Item($r, 5).Value = $eValues[$r]
}

$ws.Range("E2:E87").NumberFormat = "m/d/yy h:mm"

# ---------------------------------------------------------------------------
# 2. Column P: new "UnixDateTime" column.
# ---------------------------------------------------------------------------
$ws.Range("P1").Value = "UnixDateTime"

$ws.Range("P2").Formula = "= (E2 * 86400) - 2209075200"
$ws.Range("P3:P66").Formula = "= (E3 * 86400) - 2209075200"
$ws.Range("P67:P87").Formula = "= (E67 * 86400) - 2209075200"

# The formula result is a plain number (Unix epoch seconds) - make sure it
# does not inherit the date/time display format from its E-column precedent.
$ws.Range("P2:P87").Style = "Normal"

# ---------------------------------------------------------------------------
# 3. Column widths for the (now wider) Date/Time column and the new column.
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 29.3
$ws.Columns.Item(16).ColumnWidth = 15.75

# ---------------------------------------------------------------------------
# 4. Selection / window bookkeeping to mirror the authored workbook state.
# ---------------------------------------------------------------------------
$ws.Range("P2").Select()

$win = $excel.ActiveWindow
$win.Left = 1560
$win.Top = 1180
